# Apply updated "F" column (number of comments/likes, etc.) values for rows 3, 8, 9, 10
# on both the "展览" sheet and the "全部类型" sheet (which mirrors the same rows).
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 512
    8  = 3195
    9  = 4199
    10 = 106
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
